$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 836196.9
$ws.Range("I62").Value = 2063311.8
$ws.Range("K62").Value = 2063311.8
$ws.Range("M62").Value = -2062687.8

$ws.Range("H65").Value = 836196.9
$ws.Range("I65").Value = 2063311.8
$ws.Range("K65").Value = 10316559
$ws.Range("M65").Value = -10313439

$ws.Range("H80").Value = 7576016.5
$ws.Range("J80").Value = 278
$ws.Range("L80").Value = 834
$ws.Range("N80").Value = -2830

$ws.Range("H83").Value = 7576016.5
$ws.Range("J83").Value = 278
$ws.Range("L83").Value = 2502
$ws.Range("N83").Value = -12486

$ws.Range("H86").Value = 10557209
$ws.Range("J86").Value = 13370377
$ws.Range("L86").Value = 13370377
$ws.Range("N86").Value = -13372623

$ws.Range("H89").Value = 10557209
$ws.Range("J89").Value = 13370377
$ws.Range("L89").Value = 66851885
$ws.Range("N89").Value = -66863117

$ws.Range("H94").Value = 1999.5
$ws.Range("I94").Value = 1999.5
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 1999.5
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -1548.5
$ws.Range("N94").ClearContents()

$ws.Range("H112").Value = 38863.223
$ws.Range("J112").Value = 1783.8422
$ws.Range("L112").Value = 5351.5266
$ws.Range("N112").Value = -7567.5266

$ws.Range("H129").Value = 2352.25
$ws.Range("I129").Value = 2674.5
$ws.Range("J129").Value = 2030
$ws.Range("K129").Value = 8023.5
$ws.Range("L129").Value = 6090
$ws.Range("M129").Value = -3023.5
$ws.Range("N129").Value = -16090

$ws.Range("H138").Value = 4739.76
$ws.Range("I138").Value = 3636.1
$ws.Range("J138").Value = 5475.533
$ws.Range("K138").Value = 10908.3
$ws.Range("L138").Value = 16426.599
$ws.Range("M138").Value = -5768.299999999999
$ws.Range("N138").Value = -26706.599

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 360618.75
$ws.Range("I32").Value = 360618.75
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 360618.75
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -360331.75
$ws.Range("N32").ClearContents()

$ws.Range("H61").Value = 2400.5334
$ws.Range("I61").Value = 1609.4348
$ws.Range("J61").Value = 4999.857
$ws.Range("K61").Value = 1609.4348
$ws.Range("L61").Value = 4999.857
$ws.Range("M61").Value = -1397.4348
$ws.Range("N61").Value = -5423.857

$ws.Range("H110").Value = 47627264
$ws.Range("I110").Value = 50001130
$ws.Range("J110").Value = 150000
$ws.Range("K110").Value = 50001130
$ws.Range("L110").Value = 150000
$ws.Range("M110").Value = -49999085
$ws.Range("N110").Value = -154090

$ws.Range("H122").Value = 15154880
$ws.Range("I122").Value = 20835898
$ws.Range("K122").Value = 62507694
$ws.Range("M122").Value = -62505244

$ws.Range("H136").Value = 2400.5334
$ws.Range("I136").Value = 1609.4348
$ws.Range("J136").Value = 4999.857
$ws.Range("K136").Value = 4828.3044
$ws.Range("L136").Value = 14999.571
$ws.Range("M136").Value = -2278.3044
$ws.Range("N136").Value = -20099.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 27778650
$ws.Range("I80").Value = 1047.3334
$ws.Range("J80").Value = 41667452
$ws.Range("K80").Value = 1047.3334
$ws.Range("L80").Value = 41667452
$ws.Range("M80").Value = -49.33339999999998
$ws.Range("N80").Value = -41669448

$ws.Range("H83").Value = 27778650
$ws.Range("I83").Value = 1047.3334
$ws.Range("J83").Value = 41667452
$ws.Range("K83").Value = 5236.666999999999
$ws.Range("L83").Value = 208337260
$ws.Range("M83").Value = -244.6669999999995
$ws.Range("N83").Value = -208347244

$ws.Range("H86").Value = 1706.5294
$ws.Range("I86").Value = 1665.9
$ws.Range("J86").Value = 1764.5714
$ws.Range("K86").Value = 1665.9
$ws.Range("L86").Value = 1764.5714
$ws.Range("M86").Value = -542.9000000000001
$ws.Range("N86").Value = -4010.5714

$ws.Range("H89").Value = 1706.5294
$ws.Range("I89").Value = 1665.9
$ws.Range("J89").Value = 1764.5714
$ws.Range("K89").Value = 8329.5
$ws.Range("L89").Value = 8822.857
$ws.Range("M89").Value = -2713.5
$ws.Range("N89").Value = -20054.857

$ws.Range("H107").Value = 55598780
$ws.Range("I107").Value = 35402.2
$ws.Range("K107").Value = 35402.2
$ws.Range("M107").Value = -33482.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1606.5264
$ws.Range("I22").Value = 835.1818
$ws.Range("J22").Value = 2667.125
$ws.Range("K22").Value = 835.1818
$ws.Range("L22").Value = 2667.125
$ws.Range("M22").Value = -485.1818
$ws.Range("N22").Value = -3367.125

$ws.Range("H62").Value = 8692.615
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 8692.615
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H134").Value = 2534.0625
$ws.Range("I134").Value = 1560.6666
$ws.Range("J134").Value = 3785.5715
$ws.Range("K134").Value = 4681.9998
$ws.Range("L134").Value = 11356.7145
$ws.Range("M134").Value = -2146.9998
$ws.Range("N134").Value = -16426.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 246
$ws.Range("J33").Value = 196
$ws.Range("L33").Value = 1176
$ws.Range("N33").Value = -1742

$ws.Range("H68").Value = 2077.875
$ws.Range("I68").Value = 1550
$ws.Range("J68").Value = 2317.818
$ws.Range("K68").Value = 4650
$ws.Range("L68").Value = 6953.454000000001
$ws.Range("M68").Value = -3839
$ws.Range("N68").Value = -8575.454000000002

$ws.Range("H71").Value = 2077.875
$ws.Range("I71").Value = 1550
$ws.Range("J71").Value = 2317.818
$ws.Range("K71").Value = 13950
$ws.Range("L71").Value = 20860.362
$ws.Range("M71").Value = -9894
$ws.Range("N71").Value = -28972.362

$ws.Range("H122").Value = 342.6316
$ws.Range("I122").Value = 275.125
$ws.Range("K122").Value = 2476.125
$ws.Range("M122").Value = -26.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 12180.4
$ws.Range("I70").Value = 12488.5
$ws.Range("J70").Value = 11828.286
$ws.Range("K70").Value = 12488.5
$ws.Range("L70").Value = 11828.286
$ws.Range("M70").Value = -12218.5
$ws.Range("N70").Value = -12368.286

$ws.Range("H73").Value = 12180.4
$ws.Range("I73").Value = 12488.5
$ws.Range("J73").Value = 11828.286
$ws.Range("K73").Value = 12488.5
$ws.Range("L73").Value = 11828.286
$ws.Range("M73").Value = -11552.5
$ws.Range("N73").Value = -13700.286

$ws.Range("H122").Value = 24076974
$ws.Range("I122").Value = 465456.9
$ws.Range("J122").Value = 55558996
$ws.Range("K122").Value = 1396370.7
$ws.Range("L122").Value = 166676988
$ws.Range("M122").Value = -1393920.7
$ws.Range("N122").Value = -166681888

$ws.Range("H132").Value = 247508.48
$ws.Range("I132").Value = 315370.16
$ws.Range("K132").Value = 946110.48
$ws.Range("M132").Value = -943580.48

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1933.3334
$ws.Range("I22").Value = 1926.25
$ws.Range("K22").Value = 1926.25
$ws.Range("M22").Value = -1631.25

$ws.Range("H27").Value = 1933.3334
$ws.Range("I27").Value = 1926.25
$ws.Range("K27").Value = 1926.25
$ws.Range("M27").Value = -1819.25

$ws.Range("H46").Value = 2298.5334
$ws.Range("I46").Value = 1978.3125
$ws.Range("J46").Value = 2664.5
$ws.Range("K46").Value = 1978.3125
$ws.Range("L46").Value = 2664.5
$ws.Range("M46").Value = -1790.3125
$ws.Range("N46").Value = -3040.5

$ws.Range("H61").Value = 228534.6
$ws.Range("I61").Value = 268687.94
$ws.Range("J61").Value = 10559.286
$ws.Range("K61").Value = 268687.94
$ws.Range("L61").Value = 10559.286
$ws.Range("M61").Value = -268485.94
$ws.Range("N61").Value = -10963.286

$ws.Range("H113").Value = 228534.6
$ws.Range("I113").Value = 268687.94
$ws.Range("J113").Value = 10559.286
$ws.Range("K113").Value = 268687.94
$ws.Range("L113").Value = 10559.286
$ws.Range("M113").Value = -266517.94
$ws.Range("N113").Value = -14899.286

$ws.Range("H132").Value = 6763.2856
$ws.Range("I132").Value = 3694.4443
$ws.Range("J132").Value = 9064.916999999999
$ws.Range("K132").Value = 11083.3329
$ws.Range("L132").Value = 27194.751
$ws.Range("M132").Value = -8553.332900000001
$ws.Range("N132").Value = -32254.751

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 16674838
$ws.Range("I81").Value = 6299
$ws.Range("J81").Value = 18190160
$ws.Range("K81").Value = 12598
$ws.Range("L81").Value = 36380320
$ws.Range("M81").Value = -11537
$ws.Range("N81").Value = -36382442

$ws.Range("H84").Value = 16674838
$ws.Range("I84").Value = 6299
$ws.Range("J84").Value = 18190160
$ws.Range("K84").Value = 62990
$ws.Range("L84").Value = 181901600
$ws.Range("M84").Value = -57686
$ws.Range("N84").Value = -181912208

$ws.Range("H132").Value = 504335.78
$ws.Range("I132").Value = 913503.4399999999
$ws.Range("J132").Value = 4241.9443
$ws.Range("K132").Value = 2740510.32
$ws.Range("L132").Value = 12725.8329
$ws.Range("M132").Value = -2737980.32
$ws.Range("N132").Value = -17785.8329

$ws.Range("H140").Value = 40771
$ws.Range("J140").Value = 40771
$ws.Range("L140").Value = 40771
$ws.Range("N140").Value = -51131
